$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (last-changed) date lives in column C. Every data row's
# value advances by one day (e.g. 45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12).
$lastRow = $ws.UsedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = [double]$current + 1
    }
}
